{"js": "// Replace the 100 arithmetic-problem cell values in the single table on the\n// page with their updated equations, matching the authoring diff exactly.\n// Mapping is positional (row-major, same order the cells appear in the\n// document) and verified against the expected \"before\" text for safety.\nconst replacements = [[\"91-70=\", \"31+36=\"], [\"43-3=\", \"80-10=\"], [\"98-35=\", \"27+35=\"], [\"18+56=\", \"22+13=\"], [\"88-64=\", \"90-74=\"], [\"76-17=\", \"93-60=\"], [\"44+28=\", \"67+20=\"], [\"73+6=\", \"3+6=\"], [\"76-71=\", \"5+58=\"], [\"61-23=\", \"65+29=\"], [\"8+78=\", \"68-42=\"], [\"63+35=\", \"52+20=\"], [\"29+51=\", \"95-22=\"], [\"27+61=\", \"4+23=\"], [\"82-24=\", \"46+31=\"], [\"94-51=\", \"4+25=\"], [\"30+12=\", \"91-15=\"], [\"0+65=\", \"96-95=\"], [\"46-31=\", \"80-8=\"], [\"96-45=\", \"78-16=\"], [\"17+51=\", \"45-11=\"], [\"55+22=\", \"10+29=\"], [\"66+26=\", \"65+8=\"], [\"83+3=\", \"2+4=\"], [\"99-7=\", \"49-31=\"], [\"29+21=\", \"58-48=\"], [\"0+26=\", \"53-41=\"], [\"10+25=\", \"49-40=\"], [\"44-5=\", \"46+29=\"], [\"64+22=\", \"73-18=\"], [\"29+49=\", \"17+23=\"], [\"86-61=\", \"71-24=\"], [\"29+0=\", \"3+17=\"], [\"14+46=\", \"79+20=\"], [\"49+37=\", \"91-42=\"], [\"92-16=\", \"48+36=\"], [\"58-31=\", \"61-36=\"], [\"58+0=\", \"41-4=\"], [\"94-22=\", \"24+11=\"], [\"25+65=\", \"71+7=\"], [\"9+79=\", \"33-2=\"], [\"33-11=\", \"68+23=\"], [\"0+99=\", \"98-33=\"], [\"70-66=\", \"77-16=\"], [\"8+76=\", \"7-3=\"], [\"32-9=\", \"84+8=\"], [\"23+37=\", \"86-32=\"], [\"30+62=\", \"53-25=\"], [\"9+9=\", \"26+27=\"], [\"56-28=\", \"48-28=\"], [\"25+28=\", \"93-30=\"], [\"48-40=\", \"35+61=\"], [\"28-11=\", \"71+0=\"], [\"1+57=\", \"63+5=\"], [\"74-17=\", \"90-28=\"], [\"23+16=\", \"60-5=\"], [\"77-25=\", \"44+50=\"], [\"67-19=\", \"17+31=\"], [\"57+19=\", \"78+2=\"], [\"83-5=\", \"59-47=\"], [\"1+93=\", \"30+18=\"], [\"69-37=\", \"26+23=\"], [\"88-47=\", \"20+74=\"], [\"51-23=\", \"72-15=\"], [\"46+12=\", \"19-3=\"], [\"28+49=\", \"19+18=\"], [\"94+1=\", \"56-26=\"], [\"56-42=\", \"38-33=\"], [\"48-10=\", \"24+49=\"], [\"20+38=\", \"24+31=\"], [\"79-2=\", \"46-9=\"], [\"76+23=\", \"34-30=\"], [\"79-44=\", \"29+36=\"], [\"79-74=\", \"65+25=\"], [\"34-17=\", \"44-6=\"], [\"29+39=\", \"71-14=\"], [\"33+5=\", \"66-21=\"], [\"3+37=\", \"39+22=\"], [\"13+34=\", \"34+44=\"], [\"52+12=\", \"86-24=\"], [\"57+15=\", \"94-83=\"], [\"40-6=\", \"49+36=\"], [\"20+32=\", \"34-2=\"], [\"84-7=\", \"31-2=\"], [\"84-34=\", \"63+27=\"], [\"38+61=\", \"80+17=\"], [\"78-24=\", \"36+9=\"], [\"11+65=\", \"76-9=\"], [\"25-19=\", \"82-33=\"], [\"15+75=\", \"61-10=\"], [\"57-47=\", \"96-90=\"], [\"9+32=\", \"53+16=\"], [\"55+25=\", \"66+15=\"], [\"90-41=\", \"38+8=\"], [\"90-13=\", \"19+13=\"], [\"15+52=\", \"24-1=\"], [\"43+37=\", \"37-17=\"], [\"95+0=\", \"39+3=\"], [\"88-83=\", \"59-30=\"], [\"66-42=\", \"1+0=\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\nlet idx = 0;\nfor (let r = 0; r < values.length; r++) {\n  for (let c = 0; c < values[r].length; c++) {\n    if (idx < replacements.length) {\n      const [before, after] = replacements[idx];\n      if (values[r][c] === before) {\n        values[r][c] = after;\n        idx++;\n        continue;\n      }\n    }\n    // Positional mapping didn't line up (defensive fallback): look the\n    // current cell's text up by value among the known replacements.\n    const match = replacements.find(([before]) => before === values[r][c]);\n    if (match) {\n      values[r][c] = match[1];\n    }\n    idx++;\n  }\n}\n\ntable.values = values;\nawait context.sync();\n", "ps1": "# Update the 100 arithmetic-problem cell values in the single table on the\n# page with their revised equations. Each entry is addressed directly by its\n# (row, column) position in the table -- the same order the cells appear in\n# the document -- and is verified against the expected prior text before the\n# new text is written, so the mapping can never cross-match another cell.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    \"31+36=\"\n    \"80-10=\"\n    \"27+35=\"\n    \"22+13=\"\n    \"90-74=\"\n    \"93-60=\"\n    \"67+20=\"\n    \"3+6=\"\n    \"5+58=\"\n    \"65+29=\"\n    \"68-42=\"\n    \"52+20=\"\n    \"95-22=\"\n    \"4+23=\"\n    \"46+31=\"\n    \"4+25=\"\n    \"91-15=\"\n    \"96-95=\"\n    \"80-8=\"\n    \"78-16=\"\n    \"45-11=\"\n    \"10+29=\"\n    \"65+8=\"\n    \"2+4=\"\n    \"49-31=\"\n    \"58-48=\"\n    \"53-41=\"\n    \"49-40=\"\n    \"46+29=\"\n    \"73-18=\"\n    \"17+23=\"\n    \"71-24=\"\n    \"3+17=\"\n    \"79+20=\"\n    \"91-42=\"\n    \"48+36=\"\n    \"61-36=\"\n    \"41-4=\"\n    \"24+11=\"\n    \"71+7=\"\n    \"33-2=\"\n    \"68+23=\"\n    \"98-33=\"\n    \"77-16=\"\n    \"7-3=\"\n    \"84+8=\"\n    \"86-32=\"\n    \"53-25=\"\n    \"26+27=\"\n    \"48-28=\"\n    \"93-30=\"\n    \"35+61=\"\n    \"71+0=\"\n    \"63+5=\"\n    \"90-28=\"\n    \"60-5=\"\n    \"44+50=\"\n    \"17+31=\"\n    \"78+2=\"\n    \"59-47=\"\n    \"30+18=\"\n    \"26+23=\"\n    \"20+74=\"\n    \"72-15=\"\n    \"19-3=\"\n    \"19+18=\"\n    \"56-26=\"\n    \"38-33=\"\n    \"24+49=\"\n    \"24+31=\"\n    \"46-9=\"\n    \"34-30=\"\n    \"29+36=\"\n    \"65+25=\"\n    \"44-6=\"\n    \"71-14=\"\n    \"66-21=\"\n    \"39+22=\"\n    \"34+44=\"\n    \"86-24=\"\n    \"94-83=\"\n    \"49+36=\"\n    \"34-2=\"\n    \"31-2=\"\n    \"63+27=\"\n    \"80+17=\"\n    \"36+9=\"\n    \"76-9=\"\n    \"82-33=\"\n    \"61-10=\"\n    \"96-90=\"\n    \"53+16=\"\n    \"66+15=\"\n    \"38+8=\"\n    \"19+13=\"\n    \"24-1=\"\n    \"37-17=\"\n    \"39+3=\"\n    \"59-30=\"\n    \"1+0=\"\n)\n$oldValues = @(\n    \"91-70=\"\n    \"43-3=\"\n    \"98-35=\"\n    \"18+56=\"\n    \"88-64=\"\n    \"76-17=\"\n    \"44+28=\"\n    \"73+6=\"\n    \"76-71=\"\n    \"61-23=\"\n    \"8+78=\"\n    \"63+35=\"\n    \"29+51=\"\n    \"27+61=\"\n    \"82-24=\"\n    \"94-51=\"\n    \"30+12=\"\n    \"0+65=\"\n    \"46-31=\"\n    \"96-45=\"\n    \"17+51=\"\n    \"55+22=\"\n    \"66+26=\"\n    \"83+3=\"\n    \"99-7=\"\n    \"29+21=\"\n    \"0+26=\"\n    \"10+25=\"\n    \"44-5=\"\n    \"64+22=\"\n    \"29+49=\"\n    \"86-61=\"\n    \"29+0=\"\n    \"14+46=\"\n    \"49+37=\"\n    \"92-16=\"\n    \"58-31=\"\n    \"58+0=\"\n    \"94-22=\"\n    \"25+65=\"\n    \"9+79=\"\n    \"33-11=\"\n    \"0+99=\"\n    \"70-66=\"\n    \"8+76=\"\n    \"32-9=\"\n    \"23+37=\"\n    \"30+62=\"\n    \"9+9=\"\n    \"56-28=\"\n    \"25+28=\"\n    \"48-40=\"\n    \"28-11=\"\n    \"1+57=\"\n    \"74-17=\"\n    \"23+16=\"\n    \"77-25=\"\n    \"67-19=\"\n    \"57+19=\"\n    \"83-5=\"\n    \"1+93=\"\n    \"69-37=\"\n    \"88-47=\"\n    \"51-23=\"\n    \"46+12=\"\n    \"28+49=\"\n    \"94+1=\"\n    \"56-42=\"\n    \"48-10=\"\n    \"20+38=\"\n    \"79-2=\"\n    \"76+23=\"\n    \"79-44=\"\n    \"79-74=\"\n    \"34-17=\"\n    \"29+39=\"\n    \"33+5=\"\n    \"3+37=\"\n    \"13+34=\"\n    \"52+12=\"\n    \"57+15=\"\n    \"40-6=\"\n    \"20+32=\"\n    \"84-7=\"\n    \"84-34=\"\n    \"38+61=\"\n    \"78-24=\"\n    \"11+65=\"\n    \"25-19=\"\n    \"15+75=\"\n    \"57-47=\"\n    \"9+32=\"\n    \"55+25=\"\n    \"90-41=\"\n    \"90-13=\"\n    \"15+52=\"\n    \"43+37=\"\n    \"95+0=\"\n    \"88-83=\"\n    \"66-42=\"\n)\n\n$cols = $t.Columns.Count\n$i = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        if ($i -ge $newValues.Count) { break }\n        $cell = $t.Cell($r, $c)\n        $rng = $cell.Range\n        $rng.MoveEnd(1, -1) | Out-Null   # drop the trailing cell-mark character\n        if ($rng.Text -eq $oldValues[$i]) {\n            $rng.Text = $newValues[$i]\n        }\n        $i++\n    }\n}\n\n"}
